# Natmi following Dr Hou advice
# Rebuild the Fgf2-Sdc1 sending/target-cluster table as a full 3x3
# cross of {ECs, FAPs, sCs} x {ECs, FAPs, sCs} (rows 2-10), replacing
# the previous 3x2 table (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one worksheet row, columns A.. T in order:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate, Ligand average
#   expression value, Ligand total expression value, Ligand derived
#   specificity (avg/total), Receptor-expressing cells, Receptor
#   detection rate, Receptor average/total expression value, Receptor
#   derived specificity (avg/total), Edge average/total expression
#   weight, Edge average/total expression derived specificity.
$rowData = @(
    ,@("ECs","Fgf2","Sdc1","ECs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,1.482117666666667,4.446353,0.1401829251394648,0.1401829251394648,0.9578027328282221,8.620224595454,0.004852923345739947,0.004852923345739944)
    ,@("ECs","Fgf2","Sdc1","FAPs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,2.365790333333333,7.097371,0.2237632116883227,0.2237632116883226,1.528866767819778,13.759800910378,0.007746347943871679,0.007746347943871675)
    ,@("ECs","Fgf2","Sdc1","sCs",2,0.6666666666666666,0.6462393333333333,1.938718,0.03461850536298827,0.03461850536298827,3,1,6.724832333333333,20.174497,0.6360538631722126,0.6360538631722126,4.345851163871777,39.112660474846,0.02201923407337665,0.02201923407337664)
    ,@("FAPs","Fgf2","Sdc1","ECs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,1.482117666666667,4.446353,0.1401829251394648,0.1401829251394648,22.71023111785933,204.392080060734,0.1150665027375478,0.1150665027375477)
    ,@("FAPs","Fgf2","Sdc1","FAPs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,2.365790333333333,7.097371,0.2237632116883227,0.2237632116883226,36.25059363014867,326.255342671338,0.1836718001474224,0.1836718001474223)
    ,@("FAPs","Fgf2","Sdc1","sCs",3,1,15.322826,45.968478,0.8208310864042159,0.8208310864042158,3,1,6.724832333333333,20.174497,0.6360538631722126,0.6360538631722126,103.0434357228407,927.390921505566,0.5220927835192457,0.5220927835192456)
    ,@("sCs","Fgf2","Sdc1","ECs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,1.482117666666667,4.446353,0.1401829251394648,0.1401829251394648,3.999328526321333,35.993956736892,0.02026349905617711,0.0202634990561771)
    ,@("sCs","Fgf2","Sdc1","FAPs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,2.365790333333333,7.097371,0.2237632116883227,0.2237632116883226,6.383820245982667,57.454382213844,0.03234506359702857,0.03234506359702856)
    ,@("sCs","Fgf2","Sdc1","sCs",3,1,2.698388,8.095164,0.1445504082327959,0.1445504082327959,3,1,6.724832333333333,20.174497,0.6360538631722126,0.6360538631722126,18.14620687027866,163.315861832508,0.09194184557959022,0.09194184557959022)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $row = $rowData[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($startRow + $i, $c).Value2 = $row[$c - 1]
    }
}
